$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.890.64"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "'1.834.81"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'310.72"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.4627"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("D8").Value = "'0.3664"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "'0.07164"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("D10").Value = "'0.8797"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'0.07872"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "'19.60"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("D13").Value = "'1.833.42"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "'5.342"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "'6.393"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").Value = "'87.98"
$ws.Range("E16").Value = "  -5.23%  "
$ws.Range("D17").Value = "'1.007"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'0.000008726"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "'26.923.40"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "'14.46"
$ws.Range("E21").Value = "  -3.11%  "
$ws.Range("D22").Value = "'5.003"
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").Value = "'10.42"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").Value = "'1.990"
$ws.Range("E24").Value = "  +5.21%  "
$ws.Range("D25").Value = "'150.86"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").Value = "'18.25"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").Value = "'1.982"
$ws.Range("E27").Value = "  -5.10%  "
$ws.Range("D28").Value = "'113.62"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("D29").Value = "'4.937"
$ws.Range("E29").Value = "  -4.77%  "
$ws.Range("D30").Value = "'0.08848"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "'3.140"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").Value = "'0.7587"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "'4.462"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").Value = "'1.128"
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("D35").Value = "'2.629"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'1.087"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").Value = "'0.01936"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "'2.929"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").Value = "'0.05136"
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("D40").Value = "'6.905"
$ws.Range("E40").Value = "  -4.16%  "
$ws.Range("D41").Value = "'0.4981"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").Value = "'0.1595"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").Value = "'8.303"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "'0.4692"
$ws.Range("E44").Value = "  -3.87%  "
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "'10.09"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("D47").Value = "'102.37"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").Value = "'1.610"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("D49").Value = "'0.06093"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("D50").Value = "'64.69"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").Value = "'36.35"
$ws.Range("E51").Value = "  -2.47%  "
